# Updates the cryptos list (Price / Volume(1h) columns, plus a couple of
# row swaps in the ranking order) to match the latest scrape.
#
# Column D holds prices as text (e.g. "71.428.34", using '.' as both a
# thousands separator and decimal point), so purely-numeric-looking values
# are written with a leading apostrophe (classic "force text" entry) and
# then the cell style is reset back to Normal so no stray number format
# sticks around - this keeps the cell as plain text without leaving an
# explicit @ number format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.428.34'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '3.876.92'
$ws.Range('E3').Value = '  -2.63%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '''603.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').Value = '''172.72'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.39%  '
$ws.Range('D7').Value = '''0.670'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.93%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '''0.750'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '''0.178'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('D11').Value = '''54.01'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').Value = '''0.0000322'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '''11.51'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.84%  '
$ws.Range('E14').Value = '  -2.58%  '
$ws.Range('D15').Value = '''21.18'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('D16').Value = '3.882.56'
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').Value = '''13.95'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('E19').Value = '  -2.26%  '
$ws.Range('D20').Value = '71.284.07'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '''440.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '''4.80'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').Value = '''94.71'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.04%  '
$ws.Range('D24').Value = '''3.31'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.43%  '
$ws.Range('D25').Value = '''13.91'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('D26').Value = '''11.86'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').Value = '''4.04'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.01%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = '''10.50'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '''8.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +11.38%  '
$ws.Range('D31').Value = '''35.25'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').Value = '''13.59'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').Value = '''48.07'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').Value = '0.0₃0999'
$ws.Range('E35').Value = '  +11.11%  '
$ws.Range('D36').Value = '''69.46'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('D37').Value = '''635.46'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('D38').Value = '''0.440'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = '''0.147'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = '''3.28'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.63%  '
$ws.Range('D43').Value = '''2.89'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +9.08%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0473'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''3.16'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +18.40%  '
$ws.Range('D46').Value = '''10.19'
$ws.Range('D46').Style = "Normal"
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '''2.93'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -12.19%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '''0.144'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('D49').Value = '2.915.36'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E50').Value = '  -3.82%  '
$ws.Range('E51').Value = '  +2.96%  '
